$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Copy the date style from B4 (style index 1, date number format) onto the
# new date cells B5:B7, then set the actual values - this reuses the
# existing date style instead of minting a new numFmt.
$ws.Range("B4").Copy()
$ws.Range("B5:B7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 5 - 07/05/2018 ("x", "x", "Create Drop ")
$ws.Range("B5").Value = 43227
$ws.Range("C5").Value = "x"
$ws.Range("D5").Value = "x"

# Row 6 - 14/05/2018 ("x", "x", "Create Drop Insert into")
$ws.Range("B6").Value = 43234
$ws.Range("C6").Value = "x"
$ws.Range("D6").Value = "x"
# Written before E5 so "Create Drop Insert into" lands on shared-string
# index 6 and "Create Drop " lands on index 7 (matches original authoring
# order).
$ws.Range("E6").Value = "Create Drop Insert into"

$ws.Range("E5").Value = "Create Drop "

# Row 7 - 16/05/2018 ("x", "x")
$ws.Range("B7").Value = 43236
$ws.Range("C7").Value = "x"
$ws.Range("D7").Value = "x"

$ws.Range("F10").Select()
